$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the record forward: update the "last values" broadcast dates.
# The four "*_goods"/"services" rows move from 2023-06-01 to 2023-10-01,
# and the record_date moves from 2023-08-08 to 2023-11-07.
$ws.Range("B2").Value = [DateTime]"2023-10-01"
$ws.Range("B3").Value = [DateTime]"2023-10-01"
$ws.Range("B4").Value = [DateTime]"2023-10-01"
$ws.Range("B5").Value = [DateTime]"2023-10-01"
$ws.Range("B6").Value = [DateTime]"2023-11-07"
